$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Dades_Meteo")

# Plain text value updates (dates, measurements with units, temperatures, etc.)
# Percentage-looking strings need special handling below to avoid Excel
# auto-converting them into numeric percentage values.

$ws.Range("E2").Value = "2026-02-08 19:48:46"
$ws.Range("G2").Value = "109 cm"
$ws.Range("I2").Value = "5.2 mm"
$ws.Range("E3").Value = "2026-02-08 19:48:48"
$ws.Range("E4").Value = "2026-02-08 19:48:51"
$ws.Range("I4").Value = "2.9 mm"
$ws.Range("J4").Value = "1002.0 hPa"
$ws.Range("E5").Value = "2026-02-08 19:48:54"
$ws.Range("G5").Value = "119 cm"
$ws.Range("E6").Value = "2026-02-08 19:48:56"
$ws.Range("J6").Value = "1001.9 hPa"
$ws.Range("O6").Value = "10.1 °C"
$ws.Range("E7").Value = "2026-02-08 19:48:59"
$ws.Range("J7").Value = "1002.2 hPa"
$ws.Range("E8").Value = "2026-02-08 19:49:01"
$ws.Range("J8").Value = "1002.2 hPa"
$ws.Range("E9").Value = "2026-02-08 19:49:04"
$ws.Range("E10").Value = "2026-02-08 19:49:07"
$ws.Range("K10").Value = "11.2 MJ/m2"
$ws.Range("E11").Value = "2026-02-08 19:49:09"
$ws.Range("O11").Value = "4.7 °C"
$ws.Range("E12").Value = "2026-02-08 19:49:12"
$ws.Range("E13").Value = "2026-02-08 19:49:15"
$ws.Range("J13").Value = "1003.6 hPa"
$ws.Range("E14").Value = "2026-02-08 19:49:17"
$ws.Range("O14").Value = "11.2 °C"
$ws.Range("E15").Value = "2026-02-08 19:49:20"
$ws.Range("I15").Value = "0.2 mm"
$ws.Range("E16").Value = "2026-02-08 19:49:22"
$ws.Range("I16").Value = "3.4 mm"
$ws.Range("E17").Value = "2026-02-08 19:49:25"
$ws.Range("E18").Value = "2026-02-08 19:49:28"
$ws.Range("J18").Value = "1002.3 hPa"
$ws.Range("E19").Value = "2026-02-08 19:49:31"
$ws.Range("I19").Value = "13.0 mm"
$ws.Range("E20").Value = "2026-02-08 19:49:33"
$ws.Range("I20").Value = "8.5 mm"
$ws.Range("E21").Value = "2026-02-08 19:49:36"
$ws.Range("J21").Value = "1003.0 hPa"
$ws.Range("E22").Value = "2026-02-08 19:49:39"
$ws.Range("E23").Value = "2026-02-08 19:49:42"
$ws.Range("I23").Value = "4.9 mm"
$ws.Range("E24").Value = "2026-02-08 19:49:45"
$ws.Range("J24").Value = "1003.5 hPa"
$ws.Range("E25").Value = "2026-02-08 19:49:47"
$ws.Range("I25").Value = "0.7 mm"
$ws.Range("E26").Value = "2026-02-08 19:49:50"
$ws.Range("J26").Value = "1001.3 hPa"
$ws.Range("O26").Value = "3.7 °C"
$ws.Range("E27").Value = "2026-02-08 19:49:53"
$ws.Range("E28").Value = "2026-02-08 19:49:56"
$ws.Range("J28").Value = "1001.9 hPa"
$ws.Range("E29").Value = "2026-02-08 19:49:58"
$ws.Range("E30").Value = "2026-02-08 19:50:01"
$ws.Range("J30").Value = "1002.3 hPa"
$ws.Range("E31").Value = "2026-02-08 19:50:04"
$ws.Range("I31").Value = "0.2 mm"
$ws.Range("J31").Value = "1001.4 hPa"
$ws.Range("N31").Value = "7.8 °C 19:21 TU"
$ws.Range("E32").Value = "2026-02-08 19:50:07"
$ws.Range("E33").Value = "2026-02-08 19:50:09"
$ws.Range("J33").Value = "1003.3 hPa"
$ws.Range("E34").Value = "2026-02-08 19:50:12"
$ws.Range("E35").Value = "2026-02-08 19:50:15"
$ws.Range("J35").Value = "1004.3 hPa"
$ws.Range("E36").Value = "2026-02-08 19:50:18"
$ws.Range("J36").Value = "1002.3 hPa"
$ws.Range("E37").Value = "2026-02-08 19:50:20"
$ws.Range("J37").Value = "1003.2 hPa"
$ws.Range("E38").Value = "2026-02-08 19:50:23"
$ws.Range("I38").Value = "3.6 mm"
$ws.Range("E39").Value = "2026-02-08 19:50:26"
$ws.Range("I39").Value = "0.8 mm"
$ws.Range("E40").Value = "2026-02-08 19:50:29"
$ws.Range("J40").Value = "1003.6 hPa"
$ws.Range("E41").Value = "2026-02-08 19:50:31"
$ws.Range("J41").Value = "1002.4 hPa"
$ws.Range("K41").Value = "12.5 MJ/m2"
$ws.Range("E42").Value = "2026-02-08 19:50:34"
$ws.Range("E43").Value = "2026-02-08 19:50:37"
$ws.Range("O43").Value = "7.0 °C"
$ws.Range("E44").Value = "2026-02-08 19:50:40"
$ws.Range("I44").Value = "2.3 mm"
$ws.Range("E45").Value = "2026-02-08 19:50:42"
$ws.Range("J45").Value = "1004.5 hPa"
$ws.Range("O45").Value = "3.5 °C"
$ws.Range("E46").Value = "2026-02-08 19:50:45"
$ws.Range("J46").Value = "1004.0 hPa"

# Percentage text cells: force Text number format before assignment so Excel
# keeps the literal string (e.g. "68%") instead of converting it to 0.68,
# then restore the original "General" style (s="3") via a format-only paste
# from an already up-to-date neighboring cell that shares that same style.
$ws.Range("H4").NumberFormat = "@"
$ws.Range("H4").Value = "68%"
$ws.Range("H11").NumberFormat = "@"
$ws.Range("H11").Value = "80%"
$ws.Range("H13").NumberFormat = "@"
$ws.Range("H13").Value = "77%"
$ws.Range("H15").NumberFormat = "@"
$ws.Range("H15").Value = "73%"
$ws.Range("H24").NumberFormat = "@"
$ws.Range("H24").Value = "81%"
$ws.Range("H31").NumberFormat = "@"
$ws.Range("H31").Value = "75%"
$ws.Range("H32").NumberFormat = "@"
$ws.Range("H32").Value = "92%"
$ws.Range("H45").NumberFormat = "@"
$ws.Range("H45").Value = "77%"
$ws.Range("H46").NumberFormat = "@"
$ws.Range("H46").Value = "73%"

$ws.Range("G4").Copy() | Out-Null
$ws.Range("H4").PasteSpecial(-4122) | Out-Null
$ws.Range("H11").PasteSpecial(-4122) | Out-Null
$ws.Range("H13").PasteSpecial(-4122) | Out-Null
$ws.Range("H15").PasteSpecial(-4122) | Out-Null
$ws.Range("H24").PasteSpecial(-4122) | Out-Null
$ws.Range("H31").PasteSpecial(-4122) | Out-Null
$ws.Range("H32").PasteSpecial(-4122) | Out-Null
$ws.Range("H45").PasteSpecial(-4122) | Out-Null
$ws.Range("H46").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

